$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 86: 2023 Jan
$ws.Cells.Item(86, 1).Value = 2023
$ws.Cells.Item(86, 2).Value = "Jan"
$ws.Cells.Item(86, 3).Value = 4.41
$ws.Cells.Item(86, 4).Value = 5.19
$ws.Cells.Item(86, 5).Value = 5.09
$ws.Cells.Item(86, 6).Value = 4.94
$ws.Cells.Item(86, 7).Value = 4.7

# Row 87: 2023 Feb
$ws.Cells.Item(87, 1).Value = 2023
$ws.Cells.Item(87, 2).Value = "Feb"
$ws.Cells.Item(87, 3).Value = 4.25
$ws.Cells.Item(87, 4).Value = 5.19
$ws.Cells.Item(87, 5).Value = 5.19
$ws.Cells.Item(87, 6).Value = 5.19
$ws.Cells.Item(87, 7).Value = 4.95

[void]$ws.Range("G88").Select()
